$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 243.33333
$ws.Range("J2").Value = 620
$ws.Range("L2").Value = 620
$ws.Range("N2").Value = -846
$ws.Range("H4").Value = 209.78572
$ws.Range("I4").Value = 110.53846
$ws.Range("K4").Value = 110.53846
$ws.Range("M4").Value = 3.461539999999999
$ws.Range("H6").Value = 101708.8
$ws.Range("I6").Value = 113008.11
$ws.Range("K6").Value = 339024.33
$ws.Range("M6").Value = -338912.33
$ws.Range("H9").Value = 444
$ws.Range("J9").Value = 1349.6666
$ws.Range("L9").Value = 1349.6666
$ws.Range("N9").Value = -1687.6666
$ws.Range("H19").Value = 4993.727
$ws.Range("I19").Value = 3933.25
$ws.Range("J19").Value = 5599.7144
$ws.Range("K19").Value = 3933.25
$ws.Range("L19").Value = 5599.7144
$ws.Range("M19").Value = -3758.25
$ws.Range("N19").Value = -5949.7144
$ws.Range("H40").Value = 30004836
$ws.Range("I40").Value = 6598.5
$ws.Range("J40").Value = 50003660
$ws.Range("K40").Value = 6598.5
$ws.Range("L40").Value = 50003660
$ws.Range("M40").Value = -6423.5
$ws.Range("N40").Value = -50004010
$ws.Range("H51").Value = 9786.186
$ws.Range("I51").Value = 14454.857
$ws.Range("K51").Value = 14454.857
$ws.Range("M51").Value = -13970.857
$ws.Range("H62").Value = 27785998
$ws.Range("I62").Value = 50007000
$ws.Range("J62").Value = 9746.25
$ws.Range("K62").Value = 50007000
$ws.Range("L62").Value = 9746.25
$ws.Range("M62").Value = -50006376
$ws.Range("N62").Value = -10994.25
$ws.Range("H64").Value = 3360.5
$ws.Range("I64").Value = 3392
$ws.Range("J64").Value = 3203
$ws.Range("K64").Value = 3392
$ws.Range("L64").Value = 3203
$ws.Range("M64").Value = -3144
$ws.Range("N64").Value = -3699
$ws.Range("H65").Value = 27785998
$ws.Range("I65").Value = 50007000
$ws.Range("J65").Value = 9746.25
$ws.Range("K65").Value = 250035000
$ws.Range("L65").Value = 48731.25
$ws.Range("M65").Value = -250031880
$ws.Range("N65").Value = -54971.25
$ws.Range("H67").Value = 3360.5
$ws.Range("I67").Value = 3392
$ws.Range("J67").Value = 3203
$ws.Range("K67").Value = 3392
$ws.Range("L67").Value = 3203
$ws.Range("M67").Value = -2534
$ws.Range("N67").Value = -4919
$ws.Range("H76").Value = 4734
$ws.Range("I76").Value = 4745
$ws.Range("J76").Value = 4723
$ws.Range("K76").Value = 4745
$ws.Range("L76").Value = 4723
$ws.Range("M76").Value = -4430
$ws.Range("N76").Value = -5353
$ws.Range("H79").Value = 4734
$ws.Range("I79").Value = 4745
$ws.Range("J79").Value = 4723
$ws.Range("K79").Value = 4745
$ws.Range("L79").Value = 4723
$ws.Range("M79").Value = -3653
$ws.Range("N79").Value = -6907
$ws.Range("H99").Value = 115079660
$ws.Range("I99").Value = 5102379
$ws.Range("J99").Value = 500000160
$ws.Range("K99").Value = 15307137
$ws.Range("L99").Value = 1500000480
$ws.Range("M99").Value = -15305639
$ws.Range("N99").Value = -1500003476
$ws.Range("H103").Value = 683.8125
$ws.Range("I103").Value = 508.1
$ws.Range("J103").Value = 976.6667
$ws.Range("K103").Value = 1524.3
$ws.Range("L103").Value = 2930.0001
$ws.Range("M103").Value = -938.3000000000002
$ws.Range("N103").Value = -4102.0001
$ws.Range("H111").Value = 905.94116
$ws.Range("I111").Value = 861.61536
$ws.Range("J111").Value = 1050
$ws.Range("K111").Value = 2584.84608
$ws.Range("L111").Value = 3150
$ws.Range("M111").Value = 482.1539199999997
$ws.Range("N111").Value = -9284
$ws.Range("H116").Value = 6908.091
$ws.Range("I116").Value = 6000
$ws.Range("K116").Value = 6000
$ws.Range("M116").Value = -2558
$ws.Range("H118").Value = 4465629
$ws.Range("I118").Value = 5103190.5
$ws.Range("J118").Value = 2700
$ws.Range("K118").Value = 15309571.5
$ws.Range("L118").Value = 8100
$ws.Range("M118").Value = -15307914.5
$ws.Range("N118").Value = -11414
$ws.Range("H135").Value = 1131.5172
$ws.Range("I135").Value = 1143.9615
$ws.Range("K135").Value = 10295.6535
$ws.Range("M135").Value = -7760.653499999999
$ws.Range("H137").Value = 50151.375
$ws.Range("I137").Value = 75456
$ws.Range("J137").Value = 7977
$ws.Range("K137").Value = 226368
$ws.Range("L137").Value = 23931
$ws.Range("M137").Value = -223818
$ws.Range("N137").Value = -29031
$ws.Range("H138").Value = 2941.8113
$ws.Range("I138").Value = 1307.8889
$ws.Range("J138").Value = 3782.1143
$ws.Range("K138").Value = 3923.6667
$ws.Range("L138").Value = 11346.3429
$ws.Range("M138").Value = 1216.3333
$ws.Range("N138").Value = -21626.3429
$ws.Range("H141").Value = 2976.7273
$ws.Range("I141").Value = 2890.1428
$ws.Range("K141").Value = 8670.428400000001
$ws.Range("M141").Value = -3490.428400000001

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1834.1666
$ws.Range("I2").Value = 1412
$ws.Range("K2").Value = 1412
$ws.Range("M2").Value = -1299
$ws.Range("H27").Value = 9747.75
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H32").Value = 8547.823
$ws.Range("I32").Value = 5563.93
$ws.Range("K32").Value = 5563.93
$ws.Range("M32").Value = -5276.93
$ws.Range("H45").Value = 2566.0386
$ws.Range("I45").Value = 2545.1738
$ws.Range("J45").Value = 2726
$ws.Range("K45").Value = 2545.1738
$ws.Range("L45").Value = 2726
$ws.Range("M45").Value = -2168.1738
$ws.Range("N45").Value = -3480
$ws.Range("H61").Value = 5552.8184
$ws.Range("I61").Value = 2135.3333
$ws.Range("K61").Value = 2135.3333
$ws.Range("M61").Value = -1923.3333
$ws.Range("H74").Value = 83593.28999999999
$ws.Range("I74").Value = 87010.39
$ws.Range("K74").Value = 87010.39
$ws.Range("M74").Value = -86136.39
$ws.Range("H77").Value = 83593.28999999999
$ws.Range("I77").Value = 87010.39
$ws.Range("K77").Value = 435051.95
$ws.Range("M77").Value = -430683.95
$ws.Range("H88").Value = 2205.4546
$ws.Range("J88").Value = 2934.2
$ws.Range("L88").Value = 2934.2
$ws.Range("N88").Value = -3746.2
$ws.Range("H91").Value = 2205.4546
$ws.Range("J91").Value = 2934.2
$ws.Range("L91").Value = 2934.2
$ws.Range("N91").Value = -5742.2
$ws.Range("H102").Value = 2554.4666
$ws.Range("I102").Value = 2029.909
$ws.Range("J102").Value = 3997
$ws.Range("K102").Value = 2029.909
$ws.Range("L102").Value = 3997
$ws.Range("M102").Value = -407.9090000000001
$ws.Range("N102").Value = -7241
$ws.Range("H109").Value = 70375
$ws.Range("J109").Value = 70375
$ws.Range("L109").Value = 70375
$ws.Range("N109").Value = -73149
$ws.Range("H110").Value = 11961.857
$ws.Range("I110").Value = 11828.444
$ws.Range("J110").Value = 12202
$ws.Range("K110").Value = 11828.444
$ws.Range("L110").Value = 12202
$ws.Range("M110").Value = -9783.444
$ws.Range("N110").Value = -16292
$ws.Range("H116").Value = 1834.1666
$ws.Range("I116").Value = 1412
$ws.Range("K116").Value = 1412
$ws.Range("M116").Value = 882
$ws.Range("H122").Value = 2011.875
$ws.Range("I122").Value = 2070.5833
$ws.Range("K122").Value = 6211.749899999999
$ws.Range("M122").Value = -3761.749899999999
$ws.Range("H132").Value = 2818.756
$ws.Range("I132").Value = 3931
$ws.Range("K132").Value = 11793
$ws.Range("M132").Value = -9263
$ws.Range("H136").Value = 5552.8184
$ws.Range("I136").Value = 2135.3333
$ws.Range("K136").Value = 6405.999899999999
$ws.Range("M136").Value = -3855.999899999999

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1834.1666
$ws.Range("I3").Value = 1412
$ws.Range("K3").Value = 1412
$ws.Range("M3").Value = -1298
$ws.Range("H9").Value = 18000
$ws.Range("J9").Value = 18000
$ws.Range("L9").Value = 18000
$ws.Range("N9").Value = -18336
$ws.Range("H86").Value = 4930.7334
$ws.Range("I86").Value = 5307.875
$ws.Range("K86").Value = 5307.875
$ws.Range("M86").Value = -4184.875
$ws.Range("H89").Value = 4930.7334
$ws.Range("I89").Value = 5307.875
$ws.Range("K89").Value = 26539.375
$ws.Range("M89").Value = -20923.375
$ws.Range("H99").Value = 5104.5713
$ws.Range("I99").Value = 5808.8335
$ws.Range("K99").Value = 5808.8335
$ws.Range("M99").Value = -4310.8335
$ws.Range("H105").Value = 2007.5652
$ws.Range("I105").Value = 1777.9474
$ws.Range("K105").Value = 1777.9474
$ws.Range("M105").Value = -30.94740000000002
$ws.Range("H107").Value = 3638.0952
$ws.Range("I107").Value = 3718.875
$ws.Range("K107").Value = 3718.875
$ws.Range("M107").Value = -1798.875
$ws.Range("H134").Value = 2615.9722
$ws.Range("I134").Value = 2290.4546
$ws.Range("J134").Value = 2759.2
$ws.Range("K134").Value = 6871.3638
$ws.Range("L134").Value = 8277.599999999999
$ws.Range("M134").Value = -4336.3638
$ws.Range("N134").Value = -13347.6

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 83007.46000000001
$ws.Range("I4").Value = 7335.3335
$ws.Range("J4").Value = 105709.1
$ws.Range("K4").Value = 7335.3335
$ws.Range("L4").Value = 105709.1
$ws.Range("M4").Value = -7223.3335
$ws.Range("N4").Value = -105933.1
$ws.Range("H5").Value = 1590.8889
$ws.Range("I5").Value = 1322
$ws.Range("J5").Value = 1927
$ws.Range("K5").Value = 1322
$ws.Range("L5").Value = 1927
$ws.Range("M5").Value = -1210
$ws.Range("N5").Value = -2151
$ws.Range("H7").Value = 331.8889
$ws.Range("I7").Value = 195.71428
$ws.Range("K7").Value = 195.71428
$ws.Range("M7").Value = -82.71428
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H39").Value = 1325.5
$ws.Range("I39").Value = 1325.5
$ws.Range("K39").Value = 1325.5
$ws.Range("M39").Value = -934.5
$ws.Range("H49").Value = 1325.5
$ws.Range("I49").Value = 1325.5
$ws.Range("K49").Value = 1325.5
$ws.Range("M49").Value = -1143.5
$ws.Range("H58").Value = 3516.7188
$ws.Range("I58").Value = 3125.92
$ws.Range("K58").Value = 3125.92
$ws.Range("M58").Value = -2922.92
$ws.Range("H86").Value = 915804.5600000001
$ws.Range("I86").Value = 1671984.1
$ws.Range("J86").Value = 8389
$ws.Range("K86").Value = 1671984.1
$ws.Range("L86").Value = 8389
$ws.Range("M86").Value = -1670861.1
$ws.Range("N86").Value = -10635
$ws.Range("H89").Value = 915804.5600000001
$ws.Range("I89").Value = 1671984.1
$ws.Range("J89").Value = 8389
$ws.Range("K89").Value = 8359920.5
$ws.Range("L89").Value = 41945
$ws.Range("M89").Value = -8354304.5
$ws.Range("N89").Value = -53177
$ws.Range("H94").Value = 862.4167
$ws.Range("I94").Value = 585.9167
$ws.Range("K94").Value = 585.9167
$ws.Range("M94").Value = -134.9167
$ws.Range("H99").Value = 411408.53
$ws.Range("I99").Value = 775485.4
$ws.Range("J99").Value = 16991.916
$ws.Range("K99").Value = 775485.4
$ws.Range("L99").Value = 16991.916
$ws.Range("M99").Value = -773987.4
$ws.Range("N99").Value = -19987.916
$ws.Range("H105").Value = 4277.7173
$ws.Range("I105").Value = 1296.6522
$ws.Range("K105").Value = 1296.6522
$ws.Range("M105").Value = 450.3478
$ws.Range("H107").Value = 4200.7856
$ws.Range("I107").Value = 714.7143
$ws.Range("J107").Value = 5943.8213
$ws.Range("K107").Value = 714.7143
$ws.Range("L107").Value = 5943.8213
$ws.Range("M107").Value = 1205.2857
$ws.Range("N107").Value = -9783.8213
$ws.Range("H122").Value = 2679.6597
$ws.Range("I122").Value = 2508.1875
$ws.Range("K122").Value = 7524.5625
$ws.Range("M122").Value = -5074.5625
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H126").Value = 411408.53
$ws.Range("I126").Value = 775485.4
$ws.Range("J126").Value = 16991.916
$ws.Range("K126").Value = 2326456.2
$ws.Range("L126").Value = 50975.74800000001
$ws.Range("M126").Value = -2323986.2
$ws.Range("N126").Value = -55915.74800000001
$ws.Range("H132").Value = 6145.5835
$ws.Range("I132").Value = 2073.7144
$ws.Range("J132").Value = 11846.2
$ws.Range("K132").Value = 6221.1432
$ws.Range("L132").Value = 35538.60000000001
$ws.Range("M132").Value = -3691.1432
$ws.Range("N132").Value = -40598.60000000001
$ws.Range("H134").Value = 3162.3333
$ws.Range("I134").Value = 3018.9666
$ws.Range("K134").Value = 9056.899800000001
$ws.Range("M134").Value = -6521.899800000001
$ws.Range("H136").Value = 3516.7188
$ws.Range("I136").Value = 3125.92
$ws.Range("K136").Value = 9377.76
$ws.Range("M136").Value = -6827.76

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1261
$ws.Range("I5").Value = 974.3333
$ws.Range("K5").Value = 2922.9999
$ws.Range("M5").Value = -2810.9999
$ws.Range("H22").Value = 1994.5
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H23").Value = 829.64703
$ws.Range("I23").Value = 31.6
$ws.Range("J23").Value = 1162.1666
$ws.Range("K23").Value = 94.80000000000001
$ws.Range("L23").Value = 3486.4998
$ws.Range("M23").Value = 140.2
$ws.Range("N23").Value = -3956.4998
$ws.Range("H27").Value = 1994.5
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H37").Value = 333413300
$ws.Range("J37").Value = 333413300
$ws.Range("L37").Value = 1000239900
$ws.Range("N37").Value = -1000240124
$ws.Range("H40").Value = 307.9091
$ws.Range("J40").Value = 331.66666
$ws.Range("L40").Value = 1326.66664
$ws.Range("N40").Value = -1464.66664
$ws.Range("H76").Value = 300182800
$ws.Range("I76").Value = 375226750
$ws.Range("K76").Value = 1125680250
$ws.Range("M76").Value = -1125679867
$ws.Range("H79").Value = 300182800
$ws.Range("I79").Value = 375226750
$ws.Range("K79").Value = 1125680250
$ws.Range("M79").Value = -1125678924
$ws.Range("H135").Value = 1261
$ws.Range("I135").Value = 974.3333
$ws.Range("K135").Value = 8768.9997
$ws.Range("M135").Value = -6233.9997
$ws.Range("H136").Value = 2797.5
$ws.Range("I136").Value = 2797.5
$ws.Range("K136").Value = 8392.5
$ws.Range("M136").Value = -3292.5
$ws.Range("H137").Value = 1239.8
$ws.Range("I137").Value = 1238
$ws.Range("K137").Value = 3714
$ws.Range("M137").Value = 1386
$ws.Range("H138").Value = 2565.3333
$ws.Range("I138").Value = 2546.95
$ws.Range("K138").Value = 7640.849999999999
$ws.Range("M138").Value = -2500.849999999999
$ws.Range("H140").Value = 3630.0908
$ws.Range("J140").Value = 4001.8823
$ws.Range("L140").Value = 12005.6469
$ws.Range("N140").Value = -22365.6469

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 270.1
$ws.Range("I2").Value = 295.125
$ws.Range("J2").Value = 170
$ws.Range("K2").Value = 295.125
$ws.Range("L2").Value = 170
$ws.Range("M2").Value = -182.125
$ws.Range("N2").Value = -396
$ws.Range("H5").Value = 19988
$ws.Range("I5").Value = 19988
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 19988
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -19876
$ws.Range("N5").ClearContents()
$ws.Range("H43").Value = 10499.5
$ws.Range("I43").Value = 10499.5
$ws.Range("K43").Value = 10499.5
$ws.Range("M43").Value = -10348.5
$ws.Range("H49").Value = 20030
$ws.Range("J49").Value = 20030
$ws.Range("L49").Value = 20030
$ws.Range("N49").Value = -20398
$ws.Range("H55").Value = 11015
$ws.Range("I55").Value = 10030
$ws.Range("J55").Value = 12000
$ws.Range("K55").Value = 10030
$ws.Range("L55").Value = 12000
$ws.Range("M55").Value = -9703
$ws.Range("N55").Value = -12654
$ws.Range("H107").Value = 31659.363
$ws.Range("I107").Value = 54140.473
$ws.Range("J107").Value = 1149.2858
$ws.Range("K107").Value = 54140.473
$ws.Range("L107").Value = 1149.2858
$ws.Range("M107").Value = -52220.473
$ws.Range("N107").Value = -4989.2858
$ws.Range("H132").Value = 38740.83
$ws.Range("I132").Value = 42702.848
$ws.Range("K132").Value = 128108.544
$ws.Range("M132").Value = -125578.544

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8041.8887
$ws.Range("I7").Value = 8041.8887
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 8041.8887
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -7929.8887
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 1194.1875
$ws.Range("I22").Value = 945.44446
$ws.Range("K22").Value = 945.44446
$ws.Range("M22").Value = -650.44446
$ws.Range("H23").Value = 6058.5
$ws.Range("I23").Value = 6742.3335
$ws.Range("J23").Value = 4007
$ws.Range("K23").Value = 6742.3335
$ws.Range("L23").Value = 4007
$ws.Range("M23").Value = -6512.3335
$ws.Range("N23").Value = -4467
$ws.Range("H27").Value = 1194.1875
$ws.Range("I27").Value = 945.44446
$ws.Range("K27").Value = 945.44446
$ws.Range("M27").Value = -838.44446
$ws.Range("H40").Value = 14917.177
$ws.Range("I40").Value = 16582.732
$ws.Range("J40").Value = 2425.5
$ws.Range("K40").Value = 16582.732
$ws.Range("L40").Value = 2425.5
$ws.Range("M40").Value = -16446.732
$ws.Range("N40").Value = -2697.5
$ws.Range("H126").Value = 8041.8887
$ws.Range("I126").Value = 8041.8887
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 24125.6661
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -21655.6661
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 5886.1333
$ws.Range("I132").Value = 5878.143
$ws.Range("K132").Value = 17634.429
$ws.Range("M132").Value = -15104.429
$ws.Range("H136").Value = 3910.5881
$ws.Range("I136").Value = 1030
$ws.Range("K136").Value = 3090
$ws.Range("M136").Value = -540

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9875.1
$ws.Range("I62").Value = 7692.75
$ws.Range("J62").Value = 11330
$ws.Range("K62").Value = 7692.75
$ws.Range("L62").Value = 11330
$ws.Range("M62").Value = -7068.75
$ws.Range("N62").Value = -12578
$ws.Range("H65").Value = 9875.1
$ws.Range("I65").Value = 7692.75
$ws.Range("J65").Value = 11330
$ws.Range("K65").Value = 38463.75
$ws.Range("L65").Value = 56650
$ws.Range("M65").Value = -35343.75
$ws.Range("N65").Value = -62890
$ws.Range("H107").Value = 414.36365
$ws.Range("I107").Value = 414.36365
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1243.09095
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 676.90905
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 1061.3846
$ws.Range("I113").Value = 1080.6
$ws.Range("K113").Value = 3241.8
$ws.Range("M113").Value = -1071.8
$ws.Range("H122").Value = 2295.0417
$ws.Range("I122").Value = 2528.842
$ws.Range("K122").Value = 7586.526
$ws.Range("M122").Value = -5136.526
$ws.Range("H126").Value = 3718.2942
$ws.Range("I126").Value = 3515.2856
$ws.Range("J126").Value = 4665.6665
$ws.Range("K126").Value = 10545.8568
$ws.Range("L126").Value = 13996.9995
$ws.Range("M126").Value = -8075.856800000001
$ws.Range("N126").Value = -18936.9995
$ws.Range("H132").Value = 1155.8823
$ws.Range("I132").Value = 1086.8182
$ws.Range("J132").Value = 1282.5
$ws.Range("K132").Value = 3260.4546
$ws.Range("L132").Value = 3847.5
$ws.Range("M132").Value = -730.4546
$ws.Range("N132").Value = -8907.5
$ws.Range("H136").Value = 913072.8
$ws.Range("I136").Value = 1113755.6
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 3341266.8
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -3338716.8
$ws.Range("N136").Value = -35100
